# Updated cryptos list values (Price and Volume(1h) columns) to match the
# refreshed data in the target workbook revision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All of these source cells are plain text (inline strings) in the original
# workbook, not numbers. Several of the new "Price" values (e.g. "39.39",
# "0.320") would otherwise be auto-converted to numeric values by Excel,
# losing formatting such as trailing zeros. To avoid that, we force the
# cell to Text format before writing, then clear the format change back off
# so the cell keeps its original (default) style, as in the target diff.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

$ws.Range("D2").Value = "34.498.52"
$ws.Range("E2").Value = "  -2.92%  "
$ws.Range("D3").Value = "1.801.58"
$ws.Range("E3").Value = "  -2.31%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("E5").Value = "  -1.52%  "
$ws.Range("E6").Value = "  -1.39%  "
Set-TextValue $ws.Range("D8") "39.39"
$ws.Range("E8").Value = "  -11.46%  "
Set-TextValue $ws.Range("D9") "0.320"
$ws.Range("E9").Value = "  +2.71%  "
Set-TextValue $ws.Range("D10") "0.0679"
$ws.Range("E10").Value = "  -2.91%  "
$ws.Range("E11").Value = "  -2.11%  "
$ws.Range("D12").Value = "2.061.08"
$ws.Range("E12").Value = "  -2.41%  "
Set-TextValue $ws.Range("D13") "11.11"
$ws.Range("E13").Value = "  -2.03%  "
Set-TextValue $ws.Range("D14") "0.659"
$ws.Range("E14").Value = "  -2.56%  "
$ws.Range("D15").Value = "1.798.35"
$ws.Range("E15").Value = "  -2.76%  "
$ws.Range("E16").Value = "  -3.77%  "
$ws.Range("D17").Value = "34.361.97"
$ws.Range("E17").Value = "  -3.28%  "
Set-TextValue $ws.Range("D18") "69.01"
$ws.Range("E18").Value = "  -2.26%  "
$ws.Range("E19").Value = "  -3.02%  "
Set-TextValue $ws.Range("D20") "239.52"
$ws.Range("E20").Value = "  -1.99%  "
Set-TextValue $ws.Range("D21") "11.79"
$ws.Range("E21").Value = "  -2.64%  "
Set-TextValue $ws.Range("D22") "4.69"
$ws.Range("E22").Value = "  +0.97%  "
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("E24").Value = "  -0.98%  "
Set-TextValue $ws.Range("D25") "173.00"
$ws.Range("E25").Value = "  +0.99%  "
Set-TextValue $ws.Range("D26") "7.74"
$ws.Range("E26").Value = "  -3.49%  "
Set-TextValue $ws.Range("D27") "17.18"
$ws.Range("E28").Value = "  -0.63%  "
$ws.Range("E29").Value = "  -5.05%  "
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("E31").Value = "  +1.19%  "
Set-TextValue $ws.Range("D32") "0.0542"
$ws.Range("E32").Value = "  -2.22%  "
$ws.Range("E33").Value = "  -5.61%  "
$ws.Range("E34").Value = "  +7.54%  "
$ws.Range("E35").Value = "  -3.19%  "
Set-TextValue $ws.Range("D36") "0.693"
$ws.Range("E36").Value = "  -0.63%  "
Set-TextValue $ws.Range("D37") "90.54"
$ws.Range("E37").Value = "  -5.76%  "
$ws.Range("E38").Value = "  +5.13%  "
$ws.Range("D39").Value = "1.323.20"
$ws.Range("E39").Value = "  -1.82%  "
$ws.Range("E40").Value = "  -3.04%  "
$ws.Range("E41").Value = "  -5.87%  "
Set-TextValue $ws.Range("D42") "14.18"
$ws.Range("E42").Value = "  -7.20%  "
$ws.Range("E43").Value = "  -3.05%  "
$ws.Range("E44").Value = "  -9.49%  "
$ws.Range("E45").Value = "  -3.66%  "
Set-TextValue $ws.Range("D46") "6.13"
$ws.Range("E46").Value = "  -2.03%  "
$ws.Range("E47").Value = "  -1.30%  "
$ws.Range("D48").Value = "1.985.08"
$ws.Range("E48").Value = "  -1.53%  "
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("E50").Value = "  +3.49%  "
Set-TextValue $ws.Range("D51") "97.55"
$ws.Range("E51").Value = "  -4.92%  "
